# Apply updated crypto price/volume data to sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "65.532.29"
$ws.Range("E2").Value = "  -0.42%  "

# Row 3
$ws.Range("D3").Value = "2.653.17"
$ws.Range("E3").Value = "  -0.96%  "

# Row 4
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.87%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.53"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.58%  "

# Row 7
$ws.Range("E7").Value = "  +0.08%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.626"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +6.62%  "

# Row 9
$ws.Range("E9").Value = "  +2.81%  "

# Row 10
$ws.Range("E10").Value = "  -0.67%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.79"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.10%  "

# Row 12
$ws.Range("E12").Value = "  +0.44%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "28.76"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.51%  "

# Row 14
$ws.Range("E14").Value = "  -3.34%  "

# Row 15
$ws.Range("D15").Value = "3.130.39"
$ws.Range("E15").Value = "  -0.90%  "

# Row 16
$ws.Range("D16").Value = "65.416.62"
$ws.Range("E16").Value = "  -0.29%  "

# Row 17
$ws.Range("D17").Value = "2.667.06"
$ws.Range("E17").Value = "  -1.19%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.62"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.01%  "

# Row 19
$ws.Range("E19").Value = "  -1.27%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.46"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.91%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "348.65"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.93%  "

# Row 22
$ws.Range("E22").Value = "  -0.02%  "

# Row 23
$ws.Range("E23").Value = "  -1.64%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000112"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.45%  "

# Row 25
$ws.Range("E25").Value = "  -2.10%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.65"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.65%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.58"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.57%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.164"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.22%  "

# Row 29
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.12%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.94"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.88%  "

# Row 31
$ws.Range("B31").Value = "Bittensor"
$ws.Range("C31").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "538.98"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.81%  "

# Row 32
$ws.Range("E32").Value = "  -2.72%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.75"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.31%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.40"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.52%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.43"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.23%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.419"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.64%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.37"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.68%  "

# Row 38
$ws.Range("E38").Value = "  -0.04%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "155.77"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.23%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.91"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.53%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.03%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "161.04"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.99%  "

# Row 43
$ws.Range("E43").Value = "  -0.86%  "

# Row 44
$ws.Range("E44").Value = "  +2.87%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0605"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.67%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.49"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.73%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.635"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.41%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0255"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.02%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0993"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.67%  "

# Row 50
$ws.Range("E50").Value = "  +6.49%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.61"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.55%  "

